$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Notes row
$ws.Cells.Item(82, 3).Value = "Notes"
$ws.Cells.Item(82, 5).Value = "Rerunning with cel level selection for all marker expression"

# Fill column E (oldSelection/newSelection) first for the new temp rows
$ws.Cells.Item(83, 5).Value = "oldSelection"
$ws.Cells.Item(84, 5).Value = "newSelection"
$ws.Cells.Item(85, 5).Value = "oldSelection"
$ws.Cells.Item(86, 5).Value = "newSelection"
$ws.Cells.Item(87, 5).Value = "oldSelection"
$ws.Cells.Item(88, 5).Value = "newSelection"
$ws.Cells.Item(89, 5).Value = "oldSelection"
$ws.Cells.Item(90, 5).Value = "newSelection"
$ws.Cells.Item(91, 5).Value = "oldSelection"
$ws.Cells.Item(92, 5).Value = "newSelection"
$ws.Cells.Item(93, 5).Value = "oldSelection"
$ws.Cells.Item(94, 5).Value = "newSelection"

# Fill column B (pipeline / result names) next
$ws.Cells.Item(83, 2).Value = "Pipe_29-06"
$ws.Cells.Item(84, 2).Value = "Pipe_29-07"
$ws.Cells.Item(85, 2).Value = "Pipe_29-08"
$ws.Cells.Item(86, 2).Value = "Pipe_29-09"
$ws.Cells.Item(87, 2).Value = "Pipe_29-09"
$ws.Cells.Item(88, 2).Value = "Pipe_29-09"
$ws.Cells.Item(91, 2).Value = "SingleR_RData_2022-07-01 13-24-47"
$ws.Cells.Item(92, 2).Value = "SingleR_RData_2022-07-01 13-25-24"
$ws.Cells.Item(93, 2).Value = "SingleR_RData_2022-07-01 13-26-42"
$ws.Cells.Item(94, 2).Value = "SingleR_RData_2022-07-01 13-27-20"

# Fill remaining columns (A, C, D) - values already exist in the shared string table
$ws.Cells.Item(83, 1).Value = "results"
$ws.Cells.Item(83, 3).Value = "integration"
$ws.Cells.Item(83, 4).Value = "A+C"

$ws.Cells.Item(84, 1).Value = "results"
$ws.Cells.Item(84, 3).Value = "integration"
$ws.Cells.Item(84, 4).Value = "A+C"

$ws.Cells.Item(85, 1).Value = "results"
$ws.Cells.Item(85, 3).Value = "integration"
$ws.Cells.Item(85, 4).Value = "N+C"

$ws.Cells.Item(86, 1).Value = "results"
$ws.Cells.Item(86, 3).Value = "integration"
$ws.Cells.Item(86, 4).Value = "N+C"

$ws.Cells.Item(87, 1).Value = "results"
$ws.Cells.Item(87, 3).Value = "DEG"
$ws.Cells.Item(87, 4).Value = "A+C"

$ws.Cells.Item(88, 1).Value = "results"
$ws.Cells.Item(88, 3).Value = "DEG"
$ws.Cells.Item(88, 4).Value = "A+C"

$ws.Cells.Item(89, 3).Value = "DEG"
$ws.Cells.Item(89, 4).Value = "N+C"

$ws.Cells.Item(90, 3).Value = "DEG"
$ws.Cells.Item(90, 4).Value = "N+C"

$ws.Cells.Item(91, 1).Value = "Kriegstein"
$ws.Cells.Item(91, 3).Value = "Kriegstein to SingleR"
$ws.Cells.Item(91, 4).Value = "A+C"

$ws.Cells.Item(92, 1).Value = "Kriegstein"
$ws.Cells.Item(92, 3).Value = "Kriegstein to SingleR"
$ws.Cells.Item(92, 4).Value = "A+C"

$ws.Cells.Item(93, 1).Value = "Kriegstein"
$ws.Cells.Item(93, 3).Value = "Kriegstein to SingleR"
$ws.Cells.Item(93, 4).Value = "N+C"

$ws.Cells.Item(94, 1).Value = "Kriegstein"
$ws.Cells.Item(94, 3).Value = "Kriegstein to SingleR"
$ws.Cells.Item(94, 4).Value = "N+C"

# Update the scrolled view / active selection to match the new bottom of the log
$ws.Application.ActiveWindow.ScrollRow = 73
$ws.Range("B95").Select()
